# Update "想去人数" (interested-count) figures to the newly scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5711
$wsExpo.Range("F3").Value = 381
$wsExpo.Range("F4").Value = 644
$wsExpo.Range("F6").Value = 856
$wsExpo.Range("F8").Value = 388

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 54

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5711
$wsAll.Range("F3").Value = 381
$wsAll.Range("F4").Value = 644
$wsAll.Range("F6").Value = 856
$wsAll.Range("F8").Value = 54
$wsAll.Range("F9").Value = 388
